# Append the 2026/02/12 "逃离鸭科夫" mod-count row (row 94) to the
# "ModCounts" sheet, matching the existing data pattern (rows 3-93):
#   A = date text, B = game name text, C = mod count number,
#   all three cells center/center aligned like the row above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 93
$newRow  = $lastRow + 1

# Column A holds a literal date-like string ("2026/02/12"), not a real
# date serial. Force text storage by pre-formatting the cell as Text
# before assigning the value - otherwise Excel's normal typed-input
# parsing would turn the slash-separated string into a date.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2026/02/12"

$ws.Range("B$newRow").Value = "逃离鸭科夫"
$ws.Range("C$newRow").Value = 1190

# Match the look of the preceding data rows: copy their formatting
# (center/center alignment, General number format) onto the new row,
# without disturbing the values just written.
$ws.Range("A$lastRow`:C$lastRow").Copy()
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122)
